# Updated directory path in MarsResource, Added my login credentials in
# TestData.xls, updated Log-in source row number from TestData, Update
# worksamples directory.

$wb = $excel.ActiveWorkbook

# --- Add new login credentials row to the "LogIn" sheet -------------------
$login = $wb.Worksheets.Item("LogIn")

$login.Range("A4").Value = "http://localhost:5000/"
$login.Range("B4").Value = "dimasuhidsheila@gmail.com"
$login.Range("C4").Value = "MarsAppTest"

# Move the selection to reflect where the author last left the cursor.
$login.Range("B8").Select()

# --- Make the "LogIn" sheet the active one on re-open ----------------------
$login.Activate()
